$wb = $excel.ActiveWorkbook
$collab = $wb.Worksheets.Item("COLLABORATORS")
$newSheet = $wb.Worksheets.Add($null, $collab)
$newSheet.Name = "TYPES"

$newSheet.Range("A1").Value = "ValidTypes"
$newSheet.Range("A2").Value = "numeric"
$newSheet.Range("A3").Value = "date"
$newSheet.Range("A4").Value = "text"
$newSheet.Range("A5").Value = "categorical"

$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:A5"), $null, 1)
$tbl.Name = "Table9"
Write-Host "Table added: " $tbl.Name
Write-Host "Display: " $tbl.DisplayName
